# Atualizei dados da bibi
# - Swap the store names "Bibi Cell Manauara" / "Bibi Cell Ponta Negra" on
#   rows 4 and 5 (A4/A5), and update the daily revenue figures for those two
#   stores plus the K-column (and recomputed row totals in AG) for rows 2, 3,
#   4, 5 and 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row labels (A4 / A5) swap -------------------------------------------
$ws.Range("A4").Value = "Bibi Cell Ponta Negra"
$ws.Range("A5").Value = "Bibi Cell Manauara"

# --- Row 2 (Bibi Cell Mundi) ----------------------------------------------
$ws.Range("K2").Value = 13073.65
$ws.Range("AG2").Value = 108865.76

# --- Row 3 (Bibi Cell Vieiralves) -----------------------------------------
$ws.Range("F3").Value = 4293
$ws.Range("K3").Value = 9197.950000000001
$ws.Range("AG3").Value = 41025.95

# --- Row 4 (now Bibi Cell Ponta Negra) ------------------------------------
$ws.Range("B4").Value = 1800.01
$ws.Range("C4").Value = 4670
$ws.Range("D4").Value = 1748.51
$ws.Range("E4").Value = 5592
$ws.Range("F4").Value = 3002
$ws.Range("G4").Value = 823
$ws.Range("H4").Value = 3138.5
$ws.Range("I4").Value = 1613
$ws.Range("J4").Value = 2786.02
$ws.Range("K4").Value = 6097.5
$ws.Range("AG4").Value = 31270.54

# --- Row 5 (now Bibi Cell Manauara) ---------------------------------------
$ws.Range("B5").Value = 3340
$ws.Range("C5").Value = 1519
$ws.Range("D5").Value = 2934
$ws.Range("E5").Value = 1819
$ws.Range("F5").Value = 2503
$ws.Range("G5").Value = 2892
$ws.Range("H5").Value = 4208.4
$ws.Range("I5").Value = 3329.9
$ws.Range("J5").Value = 4038
$ws.Range("K5").Value = 2830.9
$ws.Range("AG5").Value = 29414.2

# --- Row 6 (total) ---------------------------------------------------------
$ws.Range("F6").Value = 21643.35
$ws.Range("K6").Value = 31200
$ws.Range("AG6").Value = 210576.45
